$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new job entry (Seeking DaVinci Resolve Mentor/Trainer ...) ---

$title3 = @'
Seeking DaVinci Resolve Mentor/Trainer for YouTube/IG Video Editing - Upwork
'@

$link3 = @'
https://www.upwork.com/jobs/Seeking-DaVinci-Resolve-Mentor-Trainer-for-YouTube-Video-Editing_%7E0156e8d9b16867f12d?source=rss
'@

$desc3 = @'
I am looking for someone who can help me learn to use DaVinci Resolve in order to do some video editing tasks I need to regularly complete for two YouTube channels. 
I have some experience with DaVinci Resolve, but I've forgotten a lot, and I could use some help. 
The way in which I'd like to work together is that we meet by Zoom for 15-30 hours total over the next month or so, and we work on editing actual footage I have for two YouTube channels.
Channel 1: Mostly talking head style footage with some B-Roll
Channel 2: Mix of talking head style footage and travel vlog footage.
There are several tasks I need help with. Some of these tasks are one-off tasks, and other tasks are repeated tasks.
One off tasks
- Create a channel intro
- Create an in-video "Chapter" transition
Repeated task
- Edit out parts of the video where I misspeak
- Find and add b-roll footage/images
- Color grade
- Find and add appropriate music/sound effects (I have an Epidemic Sounds subscription)
Hourly Range
: $12.00-$36.00
Posted On
: June 14, 2024 21:18 UTC
Category
: Video Editing
Skills
:DaVinci Resolve,     Video Post-Editing    
Skills
:        DaVinci Resolve,                     Video Post-Editing            
Country
: United States
click to apply
'@
$desc3 += "`n"

$content3 = @'
I am looking for someone who can help me learn to use DaVinci Resolve in order to do some video editing tasks I need to regularly complete for two YouTube channels. <br /><br />
I have some experience with DaVinci Resolve, but I&#039;ve forgotten a lot, and I could use some help. <br /><br />
The way in which I&#039;d like to work together is that we meet by Zoom for 15-30 hours total over the next month or so, and we work on editing actual footage I have for two YouTube channels.<br /><br />
Channel 1: Mostly talking head style footage with some B-Roll<br />
Channel 2: Mix of talking head style footage and travel vlog footage.<br /><br />
There are several tasks I need help with. Some of these tasks are one-off tasks, and other tasks are repeated tasks.<br /><br />
One off tasks<br />
- Create a channel intro<br />
- Create an in-video &amp;quot;Chapter&amp;quot; transition<br /><br />
Repeated task<br />
- Edit out parts of the video where I misspeak<br />
- Find and add b-roll footage/images<br />
- Color grade<br />
- Find and add appropriate music/sound effects (I have an Epidemic Sounds subscription)<br /><br /><br /><b>Hourly Range</b>: $12.00-$36.00
<br /><b>Posted On</b>: June 14, 2024 21:18 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:DaVinci Resolve,     Video Post-Editing    
<br /><b>Skills</b>:        DaVinci Resolve,                     Video Post-Editing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Seeking-DaVinci-Resolve-Mentor-Trainer-for-YouTube-Video-Editing_%7E0156e8d9b16867f12d?source=rss">click to apply</a>
'@
$content3 += "`n"

$pubdate3 = @'
Fri, 14 Jun 2024 21:18:14 +0000
'@

$hourly3 = @'
$12.00-$36.00
'@

$postedon3 = @'
June 14, 2024 21:18 UTC
'@

$category3 = @'
Video Editing
'@

$skills3 = @'
DaVinci Resolve,     Video Post-Editing
'@

$country3 = @'
United States
'@

$ws.Range("A3").Value2 = 1
$ws.Range("B3").Value2 = $title3
$ws.Range("C3").Value2 = $link3
$ws.Hyperlinks.Add($ws.Range("C3"), $link3) | Out-Null
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Value2 = $desc3
$ws.Range("E3").Value2 = $content3
$ws.Range("F3").Value2 = $pubdate3
$ws.Range("G3").Value2 = $link3
$ws.Hyperlinks.Add($ws.Range("G3"), $link3) | Out-Null
$ws.Range("G3").Style = "Hyperlink"
$ws.Range("H3").Value2 = $hourly3
$ws.Range("J3").Value2 = $postedon3
$ws.Range("K3").Value2 = $category3
$ws.Range("L3").Value2 = $skills3
$ws.Range("M3").Value2 = $country3
$ws.Range("N3").Value2 = 1
$ws.Range("O3").Value2 = 24
$ws.Range("Q3").Value2 = 0
$ws.Range("R3").Value2 = $true

$ws.Rows.Item(3).EntireRow.AutoFit()
